$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '29.390.17'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'" + '  +0.01%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'" + '1.850.18'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + '  +0.12%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'" + '0.9996'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'" + '  +0.01%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'" + '240.33'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'" + '  +0.05%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'" + '0.6286'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'" + '  -0.08%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'" + '1.000'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'" + '  -0.01%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'" + '0.07624'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'" + '  -0.08%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = "'" + '  -1.19%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'" + '24.76'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'" + '  +1.22%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'" + '0.07741'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'" + '  -0.06%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'" + '5.036'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'" + '  +0.62%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'" + '0.6790'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'" + '  -0.04%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'" + '0.00001052'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'" + '  -3.52%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'" + '83.27'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'" + '  -0.22%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'" + '6.173'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'" + '  +0.77%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'" + '29.408.32'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'" + '  -0.05%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'" + '228.09'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'" + '  -0.19%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = "'" + '  -0.59%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D21').Value = "'" + '7.499'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'" + '  +0.73%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = "'" + '  +0.00%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'" + '158.82'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'" + '  +1.04%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'" + '0.1389'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'" + '  -0.05%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'" + '8.409'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Value = "'" + '17.69'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'" + '  +0.27%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'" + '1.400'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'" + '  +7.91%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'" + '1.463'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'" + '  -0.34%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'" + '0.05610'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'" + '  -0.33%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'" + '4.113'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'" + '  +0.09%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'" + '4.073'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'" + '  +0.77%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'" + '  +0.80%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'" + '1.836'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'" + '  -0.71%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'" + '0.7014'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'" + '  -1.10%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'" + '2.581'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'" + '  -0.27%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'" + '1.235.70'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'" + '  +0.58%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'" + '0.01807'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'" + '  +0.56%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'" + '2.717'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'" + '  -2.07%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'" + '6.387'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'" + '  -1.84%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'" + '0.9023'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'" + '  -0.72%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'" + '1.000'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'" + '  +0.04%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'" + '101.56'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'" + '  +0.10%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'" + '66.05'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'" + '  +0.08%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'" + '7.226'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'" + '  +1.20%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('B45').Value = "'" + 'TheSandbox'
$ws.Range('B45').Style = 'Normal'
$ws.Range('C45').Value = "'" + 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('C45').Style = 'Normal'
$ws.Range('D45').Value = "'" + '0.4002'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'" + '  -0.14%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('B46').Value = "'" + 'BabyDogeCoin'
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').Value = "'" + 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').Value = "'" + '0.00000000116'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'" + '  -4.88%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'" + '9.010'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'" + '  +0.01%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'" + '  -0.27%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = "'" + '  +1.04%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'" + '0.05707'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'" + '  -0.07%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'" + '0.4629'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'" + '  +0.07%  '
$ws.Range('E51').Style = 'Normal'
